$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new location rows (rows 120 and 121) following the same
# pattern as the existing "Postal Code" (ara) rows already in the sheet.

$newRows = @(
    @{ Row = 120; A = 10113; B = 10113; C = 5; D = "الرمز البريدي"; E = "BNMR"; F = "ara"; G = $true; H = "superadmin"; I = "now()" },
    @{ Row = 121; A = 10114; B = 10114; C = 5; D = "الرمز البريدي"; E = "BNMR"; F = "ara"; G = $true; H = "superadmin"; I = "now()" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}

# Update the selection to match the state after data entry: the full
# remaining rows below the newly entered data are selected.
$ws.Range("A122:XFD1048576").Select()
